# Delete the obsolete "CROANCA/Facebook/7" row (row 3). This shifts all
# subsequent rows up by one, which naturally turns the old row4..row12
# data into the new row3..row11 data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(3).Delete()

# After the shift, two of the numeric values differ from a plain shift
# (the dataset counts were updated), so correct them explicitly.
$ws.Range("D3").Value = 103
$ws.Range("D4").Value = 205
